$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old placeholder note block (merged I2:L4 note + its rows) ---
# Doing this first drops the now-orphaned rich-text shared string so the
# table re-indexes before the new header strings are appended.
$ws.Rows("2:4").Delete()

# --- Rewrite the header row text/values (order matters: it determines how
#     the shared-string table re-indexes as old entries are dropped and new
#     ones appended) ---
$ws.Range("E1").Value = "2017-M"
$ws.Range("F1").Value = "2017-F"
$ws.Range("C1").Value = "Program"
$ws.Range("D1").Value = "Major"
$ws.Range("B1").Value = "HEI Code"

# --- Apply header-row formatting to the new Program/Major columns: same
#     style as the HEI Code column (bold/white font on dark fill, wrapped
#     text) instead of the plain year-header style they had before. ---
$ws.Range("B1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Match column widths for the newly added Program/Major columns to the
#     existing HEI Code column width ---
$w = $ws.Columns("B").ColumnWidth
$ws.Columns("C").ColumnWidth = $w
$ws.Columns("D").ColumnWidth = $w

# --- Update the saved selection/active cell ---
[void]$ws.Range("J19").Select()
